# "test cases login lagout"
# Update the sample login credentials shown on Sheet1 and move the
# active selection, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the username/password sample values:
#   A2 (was Giri@gmail.com) -> pavan@gmail.com
#   B2 (was Giri@9040)      -> pavan@123
# Write B2 first, then A2, so the shared-string table and cell value
# indices land the same way Excel produced them.
$ws.Range("B2").Value = "pavan@123"
$ws.Range("A2").Value = "pavan@gmail.com"

# Move the active selection from B13 to B9
$ws.Range("B9").Select() | Out-Null
